$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log entries (row 8, 9, 10). All source cells are plain text, so force
# the cells to "Text" number format first so date-like / numeric-looking
# strings ("2024-12-12", "0", "1", "3", "2", ...) are not reinterpreted as
# dates or numbers.
$newData = @(
    @{ Row = 8;  A = "2024-12-12"; B = "test";  C = "e";     D = $null;       E = "d"; F = "0" },
    @{ Row = 9;  A = "2024-12-13"; B = "1";     C = "3";     D = "Option 3";  E = "3"; F = "0" },
    @{ Row = 10; A = "2024-12-13"; B = "sdasd"; C = "sdsds"; D = "Option 2";  E = "2"; F = "0" }
)

foreach ($entry in $newData) {
    $r = $entry.Row

    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = $entry.A

    $ws.Range("B$r").NumberFormat = "@"
    $ws.Range("B$r").Value = $entry.B

    $ws.Range("C$r").NumberFormat = "@"
    $ws.Range("C$r").Value = $entry.C

    if ($entry.D -ne $null) {
        $ws.Range("D$r").NumberFormat = "@"
        $ws.Range("D$r").Value = $entry.D
    }

    $ws.Range("E$r").NumberFormat = "@"
    $ws.Range("E$r").Value = $entry.E

    $ws.Range("F$r").NumberFormat = "@"
    $ws.Range("F$r").Value = $entry.F
}
